$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 13; existing rows 13-52 shift down to 14-53.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new data record.
$ws.Cells.Item(13, 1).Value = 3
$ws.Cells.Item(13, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(13, 3).Value = "Coquimbo"
$ws.Cells.Item(13, 4).Value = 44592
$ws.Cells.Item(13, 5).Value = 5
$ws.Cells.Item(13, 6).Value = 100112022
$ws.Cells.Item(13, 7).Value = "Arveja Verde"
$ws.Cells.Item(13, 8).Value = "Perfection"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 38
$ws.Cells.Item(13, 11).Value = 22000
$ws.Cells.Item(13, 12).Value = 22000
$ws.Cells.Item(13, 13).Value = 22000
$ws.Cells.Item(13, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(13, 15).Value = "Provincia de Talca"
$ws.Cells.Item(13, 16).Value = 880
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = "Hortaliza"
